$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 45 (pushes existing rows 45-165 down to 49-169).
# Excel's default Insert() copies formatting from the row above, which already
# carries the date number format (style) on column D, matching the rest of
# the sheet.
$ws.Rows("45:48").Insert()

# Common/template field values shared by every data row in this sheet.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100101
$producto  = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad  = "Sin especificar"
$unidad    = "$/bandeja 3 kilos"
$origen    = "Región de Arica y Parinacota"
$kgUnidad  = 3

# New weekly data for the 4 newly-inserted rows (45-48).
$newRows = @(
    @{ Row = 45; Fecha = 45259; Calidad = "Especial"; Volumen = 250; PMin = 7000; PMax = 8000; PProm = 7400; PKg = 2467 },
    @{ Row = 46; Fecha = 45259; Calidad = "Primera";  Volumen = 300; PMin = 5000; PMax = 6000; PProm = 5500; PKg = 1833 },
    @{ Row = 47; Fecha = 45259; Calidad = "Segunda";  Volumen = 250; PMin = 3000; PMax = 4000; PProm = 3400; PKg = 1133 },
    @{ Row = 48; Fecha = 45259; Calidad = "Tercera";  Volumen = 300; PMin = 2000; PMax = 3000; PProm = 2500; PKg = 833 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
